# Applies the cryptos.xlsx price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "67.129.69"
$ws.Range("E2").Value2 = "  -1.45%  "
$ws.Range("D3").Value2 = "2.472.19"
$ws.Range("E4").Value2 = "  -0.01%  "
$ws.Range("D5").Formula = "=""584.30"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value2 = "  -1.41%  "
$ws.Range("D6").Formula = "=""168.69"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value2 = "  -2.85%  "
$ws.Range("E7").Value2 = "  +0.06%  "
$ws.Range("E8").Value2 = "  -2.26%  "
$ws.Range("D9").Value2 = "2.472.09"
$ws.Range("E9").Value2 = "  -2.86%  "
$ws.Range("E10").Value2 = "  -2.93%  "
$ws.Range("E11").Value2 = "  -0.83%  "
$ws.Range("E12").Value2 = "  -2.27%  "
$ws.Range("D13").Formula = "=""0.331"""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value2 = "  -3.63%  "
$ws.Range("E14").Value2 = "  -3.49%  "
$ws.Range("D16").Value2 = "67.177.85"
$ws.Range("E16").Value2 = "  -1.16%  "
$ws.Range("D17").Formula = "=""0.0000170"""
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value2 = "  -5.03%  "
$ws.Range("D18").Value2 = "2.465.27"
$ws.Range("E18").Value2 = "  -1.92%  "
$ws.Range("D19").Formula = "=""11.09"""
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value2 = "  -6.02%  "
$ws.Range("D20").Formula = "=""7.56"""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value2 = "  -5.87%  "
$ws.Range("D21").Formula = "=""354.22"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value2 = "  -4.37%  "
$ws.Range("D22").Formula = "=""4.04"""
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value2 = "  -2.72%  "
$ws.Range("E24").Value2 = "  -3.87%  "
$ws.Range("E25").Value2 = "  -7.20%  "
$ws.Range("E26").Value2 = "  -4.96%  "
$ws.Range("D27").Formula = "=""9.25"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value2 = "  -7.29%  "
$ws.Range("E28").Value2 = "  -57.84%  "
$ws.Range("D29").Value2 = "2.592.72"
$ws.Range("E29").Value2 = "  -0.64%  "
$ws.Range("D30").Value2 = "0.0₃0907"
$ws.Range("E30").Value2 = "  -7.00%  "
$ws.Range("D31").Formula = "=""515.86"""
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value2 = "  -4.88%  "
$ws.Range("E32").Value2 = "  -8.79%  "
$ws.Range("E33").Value2 = "  -6.47%  "
$ws.Range("E34").Value2 = "  -5.22%  "
$ws.Range("E35").Value2 = "  +0.02%  "
$ws.Range("E36").Value2 = "  -9.64%  "
$ws.Range("D37").Formula = "=""159.82"""
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value2 = "  +0.08%  "
$ws.Range("E38").Value2 = "  +0.16%  "
$ws.Range("E39").Value2 = "  -4.24%  "
$ws.Range("E40").Value2 = "  -6.45%  "
$ws.Range("E41").Value2 = "  -0.05%  "
$ws.Range("E42").Value2 = "  -6.38%  "
$ws.Range("D43").Formula = "=""1.68"""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value2 = "  -6.22%  "
$ws.Range("E44").Value2 = "  -6.96%  "
$ws.Range("D45").Formula = "=""2.38"""
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value2 = "  -7.53%  "
$ws.Range("D46").Formula = "=""38.87"""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value2 = "  -1.12%  "
$ws.Range("D47").Formula = "=""140.32"""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value2 = "  -5.15%  "
$ws.Range("E49").Value2 = "  -6.87%  "
$ws.Range("E50").Value2 = "  -11.38%  "
$ws.Range("D51").Formula = "=""1.60"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value2 = "  -7.12%  "

$excel.CutCopyMode = 0

